$d = $word.ActiveDocument

# Locate the "trace evidence" sub-bullet paragraph that starts with
# "glass shards from a broken window" (it currently has no cyan
# highlighting and the document's _GoBack bookmark still sits on the
# "biological evidence" bullet above).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "glass shards from a broken window*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'glass shards from a broken window' paragraph"
}

$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range

# Rewrite the paragraph so that the paragraph mark itself (pPr/rPr) and
# both of its runs carry the same cyan highlight used by the other
# "type of evidence" bullets above it. Using InsertXML (rather than
# Range.HighlightColorIndex, which only reaches the runs) lets us also
# stamp the paragraph-mark run properties, matching a manually
# re-highlighted bullet.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="1"/>
                <w:numId w:val="1"/>
              </w:numPr>
              <w:rPr>
                <w:highlight w:val="cyan"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:highlight w:val="cyan"/>
              </w:rPr>
              <w:t>glass shards from a broken window</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:highlight w:val="cyan"/>
              </w:rPr>
              <w:t>, living room</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)

# Re-resolve the paragraph (InsertXML reseats ranges) and move the
# document's _GoBack bookmark here (it previously marked the last edit
# location on the "biological evidence" bullet). Re-adding a bookmark
# with an existing name relocates it instead of creating a duplicate.
$p2 = $d.Paragraphs.Item($targetIndex)
$startOfParagraph = $d.Range($p2.Range.Start, $p2.Range.Start)
$d.Bookmarks.Add("_GoBack", $startOfParagraph)
